$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sharedStrings-level changes: swap Fiyi / Bonaire order, and update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 16:15"

# Swap country names for rows 202/203 (Fiyi <-> Bonaire, San Eustaquio y Saba)
$a202 = $ws.Range("A202").Value()
$a203 = $ws.Range("A203").Value()
$ws.Range("A202").Value = $a203
$ws.Range("A203").Value = $a202

# --- Update numeric statistics cells per country ---
# Row 4
$ws.Range("B4").Value = 6831330
$ws.Range("C4").Value = 3029
$ws.Range("D4").Value = 4120230
$ws.Range("E4").Value = 2509665
$ws.Range("G4").Value = 87
$ws.Range("H4").Value = 201435

# Row 5
$ws.Range("B5").Value = 5141905
$ws.Range("C5").Value = 26012
$ws.Range("D5").Value = 4039986
$ws.Range("E5").Value = 1018486
$ws.Range("G5").Value = 203
$ws.Range("H5").Value = 83433

# Row 13
$ws.Range("D13").Value = 456347
$ws.Range("E13").Value = 120436
$ws.Range("G13").Value = 113
$ws.Range("H13").Value = 12229

# Row 14
$ws.Range("B14").Value = 441150
$ws.Range("C14").Value = 1863
$ws.Range("D14").Value = 413928
$ws.Range("E14").Value = 15080
$ws.Range("G14").Value = 84
$ws.Range("H14").Value = 12142

# Row 19
$ws.Range("B19").Value = 328144
$ws.Range("C19").Value = 593
$ws.Range("D19").Value = 307207
$ws.Range("E19").Value = 16538
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 4399

# Row 25
$ws.Range("B25").Value = 267442
$ws.Range("C25").Value = 577
$ws.Range("E25").Value = 18891

# Row 49
$ws.Range("B49").Value = 74987
$ws.Range("C49").Value = 224
$ws.Range("D49").Value = 72967
$ws.Range("E49").Value = 1249
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 771

# Row 59
$ws.Range("B59").Value = 49627
$ws.Range("C59").Value = 612
$ws.Range("D59").Value = 45970
$ws.Range("E59").Value = 3244
$ws.Range("G59").Value = 6
$ws.Range("H59").Value = 413

# Row 61
$ws.Range("E61").Value = 6853
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 2042

# Row 63
$ws.Range("B63").Value = 45714
$ws.Range("C63").Value = 59
$ws.Range("D63").Value = 44896
$ws.Range("E63").Value = 524

# Row 68
$ws.Range("B68").Value = 38777
$ws.Range("C68").Value = 119
$ws.Range("D68").Value = 36289
$ws.Range("E68").Value = 1917
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 571

# Row 69
$ws.Range("B69").Value = 36576
$ws.Range("C69").Value = 183
$ws.Range("D69").Value = 23611
$ws.Range("E69").Value = 12323
$ws.Range("G69").Value = 5
$ws.Range("H69").Value = 642

# Row 72
$ws.Range("B72").Value = 32695
$ws.Range("C72").Value = 82
$ws.Range("D72").Value = 31512
$ws.Range("E72").Value = 445
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 738

# Row 78
$ws.Range("B78").Value = 25822
$ws.Range("C78").Value = 886
$ws.Range("D78").Value = 13908
$ws.Range("E78").Value = 11505
$ws.Range("G78").Value = 15
$ws.Range("H78").Value = 409

# Row 79
$ws.Range("B79").Value = 24605
$ws.Range("C79").Value = 394
$ws.Range("D79").Value = 17219
$ws.Range("E79").Value = 6639
$ws.Range("G79").Value = 11
$ws.Range("H79").Value = 747

# Row 101
$ws.Range("B101").Value = 9214
$ws.Range("C101").Value = 43
$ws.Range("D101").Value = 7988
$ws.Range("E101").Value = 1153

# Row 127
$ws.Range("B127").Value = 4043
$ws.Range("C127").Value = 222
$ws.Range("D127").Value = 944
$ws.Range("E127").Value = 3053

# Row 139
$ws.Range("B139").Value = 3274
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 3043
$ws.Range("E139").Value = 218

# Row 179
$ws.Range("B179").Value = 429
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 17

# Row 202
$ws.Range("C202").Value = 4
$ws.Range("D202").Value = 17
$ws.Range("E202").Value = 14
$ws.Range("H202").Value = 1

# Row 203
$ws.Range("B203").Value = 32
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 6
$ws.Range("H203").Value = 2
